# Updates cryptos list values (Price and Volume(1h) columns) for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.235.07"
$ws.Range("E2").Value = "  +1.54%  "

$ws.Range("D3").Value = "1.655.54"
$ws.Range("E3").Value = "  +0.68%  "

$ws.Range("E4").Value = "  -0.68%  "

$ws.Range("D5").Value = "'220.10"
$ws.Range("E5").Value = "  +1.65%  "

$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("E7").Value = "  -0.62%  "

$ws.Range("E8").Value = "  +0.80%  "

$ws.Range("E9").Value = "  -0.32%  "

$ws.Range("D10").Value = "'19.62"
$ws.Range("E10").Value = "  +2.07%  "

$ws.Range("E11").Value = "  +0.36%  "

$ws.Range("D12").Value = "1.886.40"
$ws.Range("E12").Value = "  +0.66%  "

$ws.Range("D13").Value = "1.655.76"
$ws.Range("E13").Value = "  +0.95%  "

$ws.Range("D14").Value = "'4.21"
$ws.Range("E14").Value = "  +1.17%  "

$ws.Range("E15").Value = "  +0.45%  "

$ws.Range("D16").Value = "'65.87"

$ws.Range("D17").Value = "27.198.37"
$ws.Range("E17").Value = "  +1.42%  "

$ws.Range("D18").Value = "0.0₃0738"
$ws.Range("E18").Value = "  +0.40%  "

$ws.Range("D19").Value = "'220.93"
$ws.Range("E19").Value = "  +2.85%  "

$ws.Range("E20").Value = "  -0.58%  "

$ws.Range("D21").Value = "'6.73"
$ws.Range("E21").Value = "  +7.55%  "

$ws.Range("D22").Value = "'4.42"
$ws.Range("E22").Value = "  +0.75%  "

$ws.Range("D23").Value = "'2.46"
$ws.Range("E23").Value = "  -0.77%  "

$ws.Range("E24").Value = "  -0.83%  "

$ws.Range("D25").Value = "'147.92"
$ws.Range("E25").Value = "  +0.90%  "

$ws.Range("E26").Value = "  -0.52%  "

$ws.Range("D27").Value = "'7.38"
$ws.Range("E27").Value = "  +2.75%  "

$ws.Range("E28").Value = "  +0.35%  "

$ws.Range("D29").Value = "'15.99"
$ws.Range("E29").Value = "  +1.93%  "

$ws.Range("E30").Value = "  +1.33%  "

$ws.Range("E31").Value = "  +1.04%  "

$ws.Range("E32").Value = "  +0.54%  "

$ws.Range("E33").Value = "  -0.27%  "

$ws.Range("E34").Value = "  +2.28%  "

$ws.Range("D35").Value = "1.269.86"
$ws.Range("E35").Value = "  -1.25%  "

$ws.Range("D36").Value = "'2.45"
$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("D37").Value = "'0.0176"
$ws.Range("E37").Value = "  -1.34%  "

$ws.Range("E38").Value = "  +0.75%  "

$ws.Range("D39").Value = "'0.827"
$ws.Range("E39").Value = "  +0.70%  "

$ws.Range("E40").Value = "  -0.58%  "

$ws.Range("E41").Value = "  -0.02%  "

$ws.Range("E42").Value = "  +1.02%  "

$ws.Range("D43").Value = "1.799.02"
$ws.Range("E43").Value = "  +0.79%  "

$ws.Range("D44").Value = "'61.89"
$ws.Range("E44").Value = "  +0.60%  "

$ws.Range("E45").Value = "  -6.01%  "

$ws.Range("D46").Value = "'92.75"
$ws.Range("E46").Value = "  +0.98%  "

$ws.Range("D47").Value = "'1.62"
$ws.Range("E47").Value = "  +0.46%  "

$ws.Range("E48").Value = "  -0.61%  "

$ws.Range("D49").Value = "'7.65"
$ws.Range("E49").Value = "  +0.13%  "

$ws.Range("D50").Value = "'0.0977"
$ws.Range("E50").Value = "  +0.60%  "

$ws.Range("E51").Value = "  -0.03%  "
